# Atualização de bases das ligas, do dia: 28-05-2024 às 07:50
#
# The source data (columns B..AD) for several match rows was re-paired /
# re-ordered; the row "id" in column A stays attached to its row number,
# but the rest of the record (match id, teams, odds, results, etc.) swaps
# between the two rows of each pair below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is a pair of worksheet row numbers whose B:AD content
# (everything except column A) must be exchanged with one another.
$rowPairs = @(
    @(74, 75),
    @(132, 133),
    @(134, 135),
    @(140, 141),
    @(142, 143),
    @(167, 168),
    @(221, 222),
    @(230, 231)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B${r1}:AD${r1}")
    $range2 = $ws.Range("B${r2}:AD${r2}")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
